$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.170.45"
$ws.Range("E2").Value = "  +11.98%  "
$ws.Range("D3").Value = "1.819.62"
$ws.Range("E3").Value = "  +8.73%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.45"
$ws.Range("E5").Value = "  +4.01%  "
$ws.Range("E6").Value = "  +3.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +5.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.50"
$ws.Range("E9").Value = "  +7.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.282"
$ws.Range("E10").Value = "  +6.73%  "
$ws.Range("E11").Value = "  +5.50%  "
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "2.081.19"
$ws.Range("E13").Value = "  +8.54%  "
$ws.Range("D14").Value = "1.821.52"
$ws.Range("E14").Value = "  +8.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("E15").Value = "  +4.83%  "
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "34.091.36"
$ws.Range("E17").Value = "  +11.55%  "
$ws.Range("E18").Value = "  +8.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.58"
$ws.Range("E19").Value = "  +5.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "258.55"
$ws.Range("E20").Value = "  +6.09%  "
$ws.Range("D21").Value = "0.0₃0752"
$ws.Range("E21").Value = "  +4.33%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +5.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.35"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.08"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.58"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("E28").Value = "  +7.60%  "
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.89"
$ws.Range("E31").Value = "  +12.19%  "
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("E33").Value = "  +4.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.54"
$ws.Range("E34").Value = "  +7.61%  "
$ws.Range("D35").Value = "1.545.44"
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "85.26"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("E39").Value = "  +5.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.626"
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("E43").Value = "  +9.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.15"
$ws.Range("E44").Value = "  +8.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0521"
$ws.Range("E45").Value = "  +4.69%  "
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("D47").Value = "1.983.38"
$ws.Range("E47").Value = "  +9.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.73"
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.89"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.89"
$ws.Range("E51").Value = "  +22.61%  "
